# 2d_calibration_parameters.xlsx — "added peat hydro params to external file"
#
# 1. Rename the "hydro" sheet to "peat_hydro_prop".
# 2. Add a brand-new "peat" sheet at the end of the workbook holding two
#    external FiPy parameters (max_sweeps / fipy_desired_residual).
# 3. Make the new "peat" sheet the active tab (it was "channel" before).
# 4. Widen two trailing columns on "channel" that were nudged when the
#    sheet was last touched, and move its remembered selection.

$wb = $excel.ActiveWorkbook

# --- 1. rename hydro -> peat_hydro_prop -----------------------------------
$hydro = $wb.Worksheets.Item("hydro")
$hydro.Name = "peat_hydro_prop"

# --- 4. channel sheet: new column widths + selection ----------------------
$channel = $wb.Worksheets.Item("channel")
$channel.Columns.Item(13).ColumnWidth = 16.5
$channel.Columns.Item(14).ColumnWidth = 23.333333333333336
[void]$channel.Range("K13").Select()

# --- 2. add the new "peat" sheet after the last existing sheet ------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$peat = $wb.Worksheets.Add($null, $lastSheet)
$peat.Name = "peat"

$peat.Columns.Item(1).ColumnWidth = 11.5
$peat.Columns.Item(2).ColumnWidth = 18.5

$peat.Range("A1").Value = "max_sweeps"
$peat.Range("B1").Value = "fipy_desired_residual"
$peat.Range("A2").Value = 1000

# fipy_desired_residual's value is stored as literal text "0.00001" (not a
# number) in the source file, so force text formatting before typing it in,
# then drop back to General so the cell keeps the workbook's default style.
$resid = $peat.Cells.Item(2, 2)
$resid.NumberFormat = "@"
$resid.Value = "0.00001"
$resid.NumberFormat = "General"

[void]$peat.Range("D10").Select()

# --- 3. make "peat" the active tab -----------------------------------------
[void]$peat.Activate()
